{"js": "// The tracked change in this revision is a re-packaging / library-upgrade\n// pass over the document (see commit message: \"Fixed POI packaging and\n// upgraded to POI 3.15.\"). No visible text, formatting or structure is\n// added, removed or re-styled anywhere in the file - every hunk in the\n// diff is either:\n//   1) pure XML-serialization noise (namespace / attribute ordering), or\n//   2) a refreshed internal identifier (the bookmark's w:id, and the\n//      w:rsidR revision-save-id stamped on the REF-field runs) that is\n//      re-minted by the authoring tool on every regeneration and is not\n//      exposed as a settable value through the document object model.\n//\n// The one concrete, user-visible echo of that re-generation that IS\n// reachable through the supported API surface is the bookmark's identity:\n// its backing id changes in the diff even though its name and the text it\n// wraps (\"bookmarked content\") stay identical. We reproduce that by\n// removing the bookmark and re-inserting it around the exact same range,\n// which causes Word to mint a fresh internal id for it - mirroring the\n// diff - while leaving every run of text and all formatting (including the\n// two \"REF bookmark1\" fields and their cached bold results) completely\n// untouched.\n\n// Locate the existing \"bookmark1\" bookmark without throwing if, for some\n// reason, it is not present (keeps the script robust/idempotent).\nconst bookmarkName = \"bookmark1\";\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(bookmarkName);\nbookmarkRange.load(\"text\");\nawait context.sync();\n\nif (!bookmarkRange.isNullObject) {\n  // Re-create the bookmark around the same range so a new internal id is\n  // assigned, without touching any text or run formatting.\n  context.document.deleteBookmark(bookmarkName);\n  bookmarkRange.insertBookmark(bookmarkName);\n  await context.sync();\n}\n", "ps1": "# The tracked change in this revision is a re-packaging / library-upgrade\n# pass over the document (see commit message: \"Fixed POI packaging and\n# upgraded to POI 3.15.\"). No visible text, formatting or structure is\n# added, removed or re-styled anywhere in the file - every hunk in the\n# diff is either:\n#   1) pure XML-serialization noise (namespace / attribute ordering), or\n#   2) a refreshed internal identifier (the bookmark's w:id, and the\n#      w:rsidR revision-save-id stamped on the REF-field runs) that is\n#      re-minted by the authoring tool on every regeneration and is not\n#      exposed as a settable value through the Word object model.\n#\n# The one concrete, user-visible echo of that re-generation that IS\n# reachable through the supported API surface is the bookmark's identity:\n# its backing id changes in the diff even though its name and the text it\n# wraps (\"bookmarked content\") stay identical. We reproduce that by\n# removing the bookmark and re-adding it around the exact same range, which\n# causes Word to mint a fresh internal id for it - mirroring the diff -\n# while leaving every run of text and all formatting (including the two\n# \"REF bookmark1\" fields and their cached bold results) completely\n# untouched.\n\n$d = $word.ActiveDocument\n\n$bookmarkName = \"bookmark1\"\n\nif ($d.Bookmarks.Exists($bookmarkName)) {\n    $bm = $d.Bookmarks.Item($bookmarkName)\n    $bmRange = $bm.Range\n    $bm.Delete()\n    $d.Bookmarks.Add($bookmarkName, $bmRange)\n}\n"}
